$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Exception 1 (step 2)" use-case text: the bracketed clause
# changes from "[Pintor já existe no sistema]" to "[Pintor Existente]".
$ws.Range("A14").Value = "Excepção 1               (passo 2)`n[Pintor Existente]"

# The row now needs less vertical space for the shorter wrapped text.
$ws.Rows.Item(14).RowHeight = 60

# Selection moved to A15 when the file was last saved.
$ws.Range("A15").Select() | Out-Null
